# Auto-generated: corrects pinyin tone-mark annotations across paragraphs.
$d = $word.ActiveDocument

# paragraph 0
$d.Content.Find.Execute("这(zhè)几(jī)天心里(lǐ)颇不宁(níng)静。今晚在院子里(lǐ)坐着(zhuó)乘(chéng)凉(liáng)，忽然想起日日走过(guò)的(de)荷(hé)塘，在这(zhè)满月的(de)光里(lǐ)，总该另有(yǒu)一番(fān)样子吧(bā)。月亮渐(jiàn)渐(jiàn)地(de)升高了(le)，墙外马路上(shàng)孩子们的(de)欢笑，已(yǐ)经(jīng)听不见(jiàn)了(le)；妻(qī)在屋里(lǐ)拍着(zhe)闰儿，迷迷糊(hú)糊(hú)地(de)哼(hēng)着(zhuó)眠歌。我悄(qiāo)悄(qiāo)地(de)披了(le)大(dà)衫，带上(shàng)门出去。", $true, $false, $false, $false, $false, $true, 1, $false, "这(zhè)几(jǐ)天心里(li)颇不宁(níng)静。今晚在院子里(lǐ)坐着(zhe)乘(chéng)凉(liáng)，忽然想起日日走过(guò)的(de)荷(hé)塘，在这(zhè)满月的(de)光里(lǐ)，总该另有(yǒu)一番(fān)样子吧(bā)。月亮渐(jiàn)渐(jiàn)地(de)升高了(le)，墙外马路上(shang)孩子们的(de)欢笑，已(yǐ)经(jīng)听不见(jiàn)了(le)；妻(qī)在屋里(lǐ)拍着(zhe)闰儿，迷迷糊(hū)糊(hú)地(de)哼(hēng)着(zhe)眠歌。我悄(qiāo)悄(qiāo)地(dì)披了(le)大(dà)衫，带上(shang)门出去。", 2) | Out-Null

# paragraph 1
$d.Content.Find.Execute("沿(yán)着(zháo)荷(hé)塘，是一条曲(qū)折(zhé)的(de)小煤屑路。这(zhè)是一条幽僻的(de)路；白天也少(shào)人走，夜晚更(gēng)加寂寞。荷(hé)塘四面，长(cháng)着(zhe)许多树，蓊[weng3]蓊[weng3]郁郁的(de)。路的(dì)一旁，是些杨柳，和(hé)一些不知道名字的(de)树。没(méi)有(yǒu)月光的(de)晚上(shàng)，这(zhè)路上(shàng)阴森森的(de)，有(yǒu)些怕人。今晚却很好(hǎo)，虽然月光也还(huán)是淡淡的(de)。", $true, $false, $false, $false, $false, $true, 1, $false, "沿(yán)着(zhe)荷(hé)塘，是一条曲(qǔ)折(zhé)的(de)小煤屑路。这(zhè)是一条幽僻的(de)路；白天也少(shào)人走，夜晚更(gèng)加寂寞。荷(hé)塘四面，长(cháng)着(zhe)许多树，蓊[weng3]蓊[weng3]郁郁的(de)。路的(de)一旁，是些杨柳，和(hé)一些不知道名字的(de)树。没(méi)有(yǒu)月光的(de)晚上(shang)，这(zhè)路上(shang)阴森森的(de)，有(yǒu)些怕人。今晚却很好(hào)，虽然月光也还(hái)是淡淡的(de)。", 2) | Out-Null

# paragraph 2
$d.Content.Find.Execute("路上(shàng)只(zhǐ)我一个(gè)人，背(bēi)着(zhe)手踱着(zhe)。这(zhè)一片(piàn)天地(de)好(hào)像是我的(de)；我也像超出了(xi)平常的(de)自己，到了(le)另一个(gè)世界里(lǐ)。我爱热闹，也爱冷静；爱群居，也爱独处(chù)。", $true, $false, $false, $false, $false, $true, 1, $false, "路上(shang)只(zhǐ)我一个(ge)人，背(bēi)着(zháo)手踱着(zhe)。这(zhè)一片(piàn)天地(dì)好(hào)像是我的(de)；我也像超出了(le)平常的(de)自己，到了(le)另一个(ge)世界里(lǐ)。我爱热闹，也爱冷静；爱群居，也爱独处(chù)。", 2) | Out-Null

# paragraph 3
$d.Content.Find.Execute("像今晚上(shang)，一个(gè)人在这(zhè)苍茫的(de)月下，什(shí)么(me)都(dōu)可(kě)以想，什(shí)么(me)都(dū)可(kě)以不想，便(biàn)觉(jué)是个(gè)自由的(de)人。白天里(lǐ)一定要(yào)做的(de)事，一定要(yāo)说(shuō)的(de)话，现 在都(dū)可(kě)不理。这(zhè)是独处(chù)的(de)妙处(chù)，我且(qiě)受用这(zhè)无(wú)边(biān)的(de)荷(hé)香月色(sè)好(hǎo)了(liǎo)。", $true, $false, $false, $false, $false, $true, 1, $false, "像今晚上(shang)，一个(ge)人在这(zhè)苍茫的(de)月下，什(shén)么(me)都(dū)可(kě)以想，什(shén)么(me)都(dū)可(kě)以不想，便(biàn)觉(jué)是个(gè)自由的(de)人。白天里(lǐ)一定要(yāo)做的(de)事，一定要(yào)说(shuō)的(de)话，现 在都(dōu)可(kě)不理。这(zhè)是独处(chù)的(de)妙处(chù)，我且(qiě)受用这(zhè)无(wú)边(bian)的(de)荷(hé)香月色(sè)好(hào)了(le)。", 2) | Out-Null

# paragraph 4
$d.Content.Find.Execute("曲(qū)曲(qū)折(zhé)折(zhé)的(de)荷(hé)塘上(shàng)面，弥望的(de)是田田的(de)叶(yè)子。叶(yè)子出水很高，像亭亭的(dì)舞女的(dì)裙。层层的(de)叶(yè)子中(zhōng)间(jiān)，零星地(dì)点缀着(zháo)些白花，有(yǒu)袅[niao3]娜(nà)地(dì)开着(zhe)的(de)，有(yǒu)羞涩[se4]地(dì)打(dǎ)着(zhe)朵儿的(de)；正(zhèng)如一粒粒的(de)明珠，又如碧天里(lǐ)的(de)星星，又如刚出浴的(de)美人。微风过(guò)处(chù)，送来缕缕清香，仿佛(fó)远处(chù)高楼上(shàng)渺[miao3]茫的(de)歌声似(sì)的(de)。", $true, $false, $false, $false, $false, $true, 1, $false, "曲(qǔ)曲(qū)折(zhé)折(zhé)的(de)荷(hé)塘上(shàng)面，弥望的(de)是田田的(de)叶(yè)子。叶(yè)子出水很高，像亭亭的(de)舞女的(de)裙。层层的(de)叶(yè)子中(zhōng)间(jiān)，零星地(de)点缀着(zhe)些白花，有(yǒu)袅[niao3]娜(nà)地(dì)开着(zhe)的(de)，有(yǒu)羞涩[se4]地(de)打(dǎ)着(zhe)朵儿的(de)；正(zhēng)如一粒粒的(de)明珠，又如碧天里(lǐ)的(de)星星，又如刚出浴的(de)美人。微风过(guò)处(chù)，送来缕缕清香，仿佛(fú)远处(chù)高楼上(shang)渺[miao3]茫的(de)歌声似(sì)的(dì)。", 2) | Out-Null

# paragraph 5
$d.Content.Find.Execute("忽然想起采(cǎi)莲的(de)事情来了(liǎo)。采(cǎi)莲是江南(nán)的(de)旧俗，似(sì)乎很早就有(yǒu)，而六(liù)朝(zhāo)时为(wéi)盛(shèng)；从(cóng)诗歌里(lǐ)可(kě)以约(yuē)略知道。采(cǎi)莲的(de)是少(shǎo)年的(de)女子，她们是荡着(zháo)小船，唱着(zhuó)艳歌去的(de)。采(cǎi)莲人不用说(shuō)很多，还(huán)有(yǒu)看(kàn)采(cǎi)莲的(de)人。那(nà)是一个(gè)热闹的(de)季节(jié)，也是一个(gè)风流的(de)季节(jié)。梁元帝《采(cǎi)莲赋》里(lǐ)说(shuō)得(dé)好(hào)：", $true, $false, $false, $false, $false, $true, 1, $false, "忽然想起采(cǎi)莲的(de)事情来了(le)。采(cǎi)莲是江南(nán)的(de)旧俗，似(sì)乎很早就有(yǒu)，而六(lù)朝(cháo)时为(wèi)盛(shèng)；从(cóng)诗歌里(lǐ)可(kě)以约(yuē)略知道。采(cǎi)莲的(de)是少(shào)年的(de)女子，她们是荡着(zhe)小船，唱着(zhe)艳歌去的(de)。采(cǎi)莲人不用说(shuō)很多，还(huán)有(yǒu)看(kàn)采(cǎi)莲的(de)人。那(nà)是一个(ge)热闹的(de)季节(jié)，也是一个(ge)风流的(de)季节(jié)。梁元帝《采(cǎi)莲赋》里(lǐ)说(shuō)得(dé)好(hào)：", 2) | Out-Null

# paragraph 6
$d.Content.Find.Execute("于是妖童媛[yuan2](yuán)女，荡舟心许；鷁首徐回，兼传(chuán)羽杯；棹[zhuo1]将(jiāng)移而藻挂，船欲动而萍开。尔其(qí)纤(xiān)腰束素，迁延顾步；夏始春余，叶(yè)嫩花初，恐沾裳而浅(qiǎn)笑，畏倾船而敛裾[ju1]。", $true, $false, $false, $false, $false, $true, 1, $false, "于是妖童媛[yuan2](yuán)女，荡舟心许；鷁首徐回，兼传(zhuàn)羽杯；棹[zhuo1]将(jiāng)移而藻挂，船欲动而萍开。尔其(qí)纤(xiān)腰束素，迁延顾步；夏始春余，叶(yè)嫩花初，恐沾裳而浅(qiǎn)笑，畏倾船而敛裾[ju1]。", 2) | Out-Null

# paragraph 7
$d.Content.Find.Execute("可(kě)见(jiàn)当(dàng)时嬉[xi1]游的(de)光景了(le)。这(zhè)真是有(yǒu)趣(qù)的(de)事，可(kě)惜我们现在早已(yǐ)无(wú)福消受了(le)。", $true, $false, $false, $false, $false, $true, 1, $false, "可(kě)见(jiàn)当(dāng)时嬉[xi1]游的(de)光景了(le)。这(zhè)真是有(yǒu)趣(qù)的(de)事，可(kě)惜我们现在早已(yǐ)无(wú)福消受了(le)。", 2) | Out-Null

# paragraph 8
$d.Content.Find.Execute("他勤奋好(hào)学，总是取得(dé)好(hǎo)成绩。", $true, $false, $false, $false, $false, $true, 1, $false, "他勤奋好(hào)学，总是取得(dé)好(hào)成绩。", 2) | Out-Null

# paragraph 9
$d.Content.Find.Execute("他怀着(zhuó)沉重(zhòng)的(de)心情，重(chóng)新出发(fā)。", $true, $false, $false, $false, $false, $true, 1, $false, "他怀着(zhe)沉重(zhòng)的(de)心情，重(chóng)新出发(fā)。", 2) | Out-Null

# paragraph 10
$d.Content.Find.Execute("干(gàn)一行(xíng)爱一行(xíng)，你一定行(xíng)的(de)。", $true, $false, $false, $false, $false, $true, 1, $false, "干(gàn)一行(xíng)爱一行(xíng)，你一定行(háng)的(de)。", 2) | Out-Null

# paragraph 11
$d.Content.Find.Execute("知道了(le)，我了(liǎo)解(jiě)这(zhè)个(gè)事情的(de)严重(zhòng)性。", $true, $false, $false, $false, $false, $true, 1, $false, "知道了(le)，我了(liǎo)解(jiě)这(zhè)个(ge)事情的(de)严重(zhòng)性。", 2) | Out-Null
